$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 282 (shifts existing rows 282.. down to 285..)
$ws.Rows.Item(282).Resize(3).Insert()

# Populate the 3 newly inserted rows (now rows 282-284) with the new data block
# Columns: A Mercado ID, B Mercado, C Region, D Fecha, E ?, F Producto ID, G Producto,
#          H Variedad, I Categoria, J Volumen, K Precio Min, L Precio Max, M Precio Prom,
#          N Unidad, O Zona, P Precio/Kg (?), Q Kilos, R Tipo

$ws.Range("A282:A284").Value = 8
$ws.Range("B282:B284").Value = "Terminal La Palmera de La Serena"
$ws.Range("C282:C284").Value = "Coquimbo"
$ws.Range("D282:D284").Value = 44641
$ws.Range("E282:E284").Value = 4
$ws.Range("F282:F284").Value = 100112043
$ws.Range("G282:G284").Value = "Pepino dulce"
$ws.Range("H282:H284").Value = "Cultivar IV Región"
$ws.Range("N282:N284").Value = "`$/bandeja 18 kilos"
$ws.Range("O282:O284").Value = "Provincia de Limarí"
$ws.Range("Q282:Q284").Value = 18
$ws.Range("R282:R284").Value = "Hortaliza"

$ws.Range("I282").Value = "Primera"
$ws.Range("J282").Value = 600
$ws.Range("K282").Value = 9000
$ws.Range("L282").Value = 10000
$ws.Range("M282").Value = 9500
$ws.Range("P282").Value = 528

$ws.Range("I283").Value = "Segunda"
$ws.Range("J283").Value = 400
$ws.Range("K283").Value = 7000
$ws.Range("L283").Value = 8000
$ws.Range("M283").Value = 7500
$ws.Range("P283").Value = 417

$ws.Range("I284").Value = "Tercera"
$ws.Range("J284").Value = 200
$ws.Range("K284").Value = 5000
$ws.Range("L284").Value = 6000
$ws.Range("M284").Value = 5500
$ws.Range("P284").Value = 306
